$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4 and 5 repeat the same Transaction Number value already present in
# A2/A3. Copy those cells down instead of re-typing the literal so Excel
# keeps storing it as a shared-string reference (matching A2/A3) rather than
# re-interpreting the all-digit text as a number.
$ws.Range("A2:A3").Copy($ws.Range("A4:A5"))
